# Add session runs for Modal Model + KFOLD
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Modal Model (with KFolds)")

# Raw percentage-accuracy numbers for each of the 10 test folds (rows 3-12)
$accuracy = @(
    @(94.17, 94.8,  94.94, 94.75, 94.57, 95.52, 94,    95.38, 94.9,  94.64),
    @(94.38, 94.91, 95.18, 95.11, 94.82, 94.46, 94.77, 94.04, 94.98, 94.64),
    @(94.4,  94.33, 95.66, 94.47, 95.29, 94.84, 94.61, 95.6,  94.69, 94.27),
    @(93.8,  94.02, 94.64, 94.46, 94.21, 95.46, 94.35, 95.23, 94.98, 94.28),
    @(94.04, 94.9,  94.78, 94.5,  95.14, 95.2,  93.52, 95.27, 94.68, 94.67),
    @(94.11, 94.5,  94.51, 94.09, 94.52, 95.09, 94.3,  94.68, 94.05, 94.5),
    @(94.06, 94.64, 94.56, 94.69, 94.45, 94.82, 94.66, 94.03, 94.91, 94.98),
    @(93.9,  94.29, 95.2,  94.88, 94.34, 94.52, 94.65, 95.57, 94.59, 94.77),
    @(93.75, 94.54, 95.06, 94.67, 94.66, 94.86, 93.36, 95.08, 95.51, 95.1),
    @(93.75, 94.54, 95.06, 94.67, 94.66, 94.86, 94.36, 95.08, 95.51, 95.1)
)

# Confusion-matrix counts for the KFold runs (rows 13-16)
$tp = @(10, 15, 13, 14, 14, 10, 12, 12, 15, 9)
$tn = @(46, 38, 38, 37, 33, 40, 35, 37, 33, 38)
$fn = @(11, 10, 12, 14, 13, 14, 17, 14, 17, 13)
$fp = @(10, 14, 14, 12, 17, 13, 13, 14, 12, 17)

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

for ($r = 0; $r -lt 10; $r++) {
    $rowNum = $r + 3
    for ($c = 0; $c -lt 10; $c++) {
        $ws.Range($cols[$c] + $rowNum).Value = $accuracy[$r][$c]
    }
}

for ($c = 0; $c -lt 10; $c++) {
    $ws.Range($cols[$c] + "13").Value = $tp[$c]
    $ws.Range($cols[$c] + "14").Value = $tn[$c]
    $ws.Range($cols[$c] + "15").Value = $fn[$c]
    $ws.Range($cols[$c] + "16").Value = $fp[$c]
}

# Move the active tab / selection from "Mean Model (with K-Folds)" to
# "Modal Model (with KFolds)".
$meanWs = $wb.Worksheets.Item("Mean Model (with K-Folds)")
$meanWs.Range("E28").Select()
$ws.Activate()
$ws.Range("K13").Select()
